$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Guild name corrections / re-ranking ---
# "Inferno" overtook "Fabíola" in score -> their rows swap places (row5 <-> row6).
$ws.Range("A5").Value = "Inferno"
$ws.Range("A6").Value = "Fabíola"

# "supercell" and "Sergas" swap places too (row32 <-> row33).
$ws.Range("A32").Value = "supercell"
$ws.Range("A33").Value = "Sergas"

# Simple rename: "RubenVski" -> "diogo" (no row reorder for this one).
$ws.Range("A49").Value = "diogo"

# --- Row 5 (now "Inferno") stats ---
$ws.Range("B5").Value = 39.0
$ws.Range("C5").Value = 100.0
$ws.Range("D5").Value = 117.0
$ws.Range("E5").Value = 45372.0
$ws.Range("F5").Value = 38.0
$ws.Range("G5").Value = 28.0
$ws.Range("H5").Value = 0.0
$ws.Range("I5").Value = 1.0
$ws.Range("J5").Value = 74.0
$ws.Range("K5").Value = 295.0

# --- Row 6 (now "Fabíola") stats ---
$ws.Range("B6").Value = 38.0
$ws.Range("C6").Value = 100.0
$ws.Range("D6").Value = 114.0
$ws.Range("E6").Value = 43575.0
$ws.Range("F6").Value = 39.0
$ws.Range("G6").Value = 26.0
$ws.Range("H6").Value = 0.0
$ws.Range("I6").Value = 0.0
$ws.Range("J6").Value = 67.0
$ws.Range("K6").Value = 294.0

# --- WARSCORE (column K) bumps/adjustments for several guilds ---
$ws.Range("K13").Value = 259.0
$ws.Range("K19").Value = 212.0
$ws.Range("K22").Value = 192.0
$ws.Range("K24").Value = 179.0
$ws.Range("K28").Value = 154.0
$ws.Range("K30").Value = 129.0

# --- Row 32 (now "supercell") stats ---
$ws.Range("B32").Value = 25.0
$ws.Range("C32").Value = 87.0
$ws.Range("D32").Value = 73.0
$ws.Range("E32").Value = 27916.0
$ws.Range("F32").Value = 26.0
$ws.Range("G32").Value = 11.0
$ws.Range("H32").Value = 2.0
$ws.Range("I32").Value = 1.0
$ws.Range("J32").Value = 42.0
$ws.Range("K32").Value = 111.0

# --- Row 33 (now "Sergas") stats ---
$ws.Range("B33").Value = 28.0
$ws.Range("C33").Value = 100.0
$ws.Range("D33").Value = 75.0
$ws.Range("E33").Value = 30992.0
$ws.Range("F33").Value = 28.0
$ws.Range("G33").Value = 12.0
$ws.Range("H33").Value = 9.0
$ws.Range("I33").Value = 1.0
$ws.Range("J33").Value = 43.0
$ws.Range("K33").Value = 111.0

$ws.Range("K39").Value = 57.0

# --- Row 49 ("diogo") dropped out of the war, stats reset to zero/blank ---
$ws.Range("B49").Value = 0.0
$ws.Range("C49").Value = 0.0
$ws.Range("D49").Value = 0.0
$ws.Range("E49").Value = 0.0
$ws.Range("F49").Value = 0.0
$ws.Range("G49").Value = 0.0
$ws.Range("H49").Value = 0.0
$ws.Range("I49").Value = 0.0
$ws.Range("J49").ClearContents()
$ws.Range("K49").Value = 1.0

$ws.Range("K50").Value = -18.0
